$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Price 6 -> 10
$ws.Range("G5").Value = 10

# Row 6: fix the supplier import - correct item id, description, price & category
$ws.Range("E6").Value = 22
$ws.Range("F6").Value = "Exercise Book A4 Hardcover(150pg) "
$ws.Range("G6").Value = 10
$ws.Range("I6").Value = "Exercise"

# Widen column F so the longer description fits
$ws.Columns.Item(6).ColumnWidth = 47.6

# Move the active selection to G5
$ws.Range("G5").Select()
